$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows for 2025-03-26 (original rows 2-4); remaining rows shift up.
$ws.Range("A2:A4").EntireRow.Delete()

# "Atualização automática via cronjob": refresh the quantidade_atipica (A)
# and estoque_atualizado (G) columns for the remaining records.
$ws.Range("A2").Value = 9
$ws.Range("G2").Value = 331

$ws.Range("A3").Value = 4
$ws.Range("G3").Value = 40

$ws.Range("A4").Value = 8
$ws.Range("G4").Value = 352

$ws.Range("A5").Value = 0
$ws.Range("G5").Value = 201

$ws.Range("A6").Value = 1
$ws.Range("G6").Value = 1632

$ws.Range("A7").Value = 2
$ws.Range("G7").Value = 21

$ws.Range("A8").Value = 3
$ws.Range("G8").Value = 35

$ws.Range("A9").Value = 5
$ws.Range("G9").Value = 805

$ws.Range("A10").Value = 6
$ws.Range("G10").Value = 67

$ws.Range("A11").Value = 10
$ws.Range("G11").Value = 1

$ws.Range("A12").Value = 7
$ws.Range("G12").Value = 10
